$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each row 3-37, copy the formatting (style) of column P into the new
# column Q, then set Q's value to match the newly added 2020 data column.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$ws.Range("Q3").Value = 2020
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 0.1
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 0.1
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 0.1
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = 0
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 0
$ws.Range("P9").Copy()
$ws.Range("Q9").PasteSpecial(-4122)
$ws.Range("Q9").Value = 0
$ws.Range("P10").Copy()
$ws.Range("Q10").PasteSpecial(-4122)
$ws.Range("Q10").Value = 0
$ws.Range("P11").Copy()
$ws.Range("Q11").PasteSpecial(-4122)
$ws.Range("Q11").Value = 0
$ws.Range("P12").Copy()
$ws.Range("Q12").PasteSpecial(-4122)
$ws.Range("Q12").Value = 0
$ws.Range("P13").Copy()
$ws.Range("Q13").PasteSpecial(-4122)
$ws.Range("Q13").Value = 0
$ws.Range("P14").Copy()
$ws.Range("Q14").PasteSpecial(-4122)
$ws.Range("Q14").Value = 0.1
$ws.Range("P15").Copy()
$ws.Range("Q15").PasteSpecial(-4122)
$ws.Range("Q15").Value = 0
$ws.Range("P16").Copy()
$ws.Range("Q16").PasteSpecial(-4122)
$ws.Range("Q16").Value = 0
$ws.Range("P17").Copy()
$ws.Range("Q17").PasteSpecial(-4122)
$ws.Range("Q17").Value = 0
$ws.Range("P18").Copy()
$ws.Range("Q18").PasteSpecial(-4122)
$ws.Range("Q18").Value = 0
$ws.Range("P19").Copy()
$ws.Range("Q19").PasteSpecial(-4122)
$ws.Range("Q19").Value = 0
$ws.Range("P20").Copy()
$ws.Range("Q20").PasteSpecial(-4122)
$ws.Range("Q20").Value = 0
$ws.Range("P21").Copy()
$ws.Range("Q21").PasteSpecial(-4122)
$ws.Range("Q21").Value = 0
$ws.Range("P22").Copy()
$ws.Range("Q22").PasteSpecial(-4122)
$ws.Range("Q22").Value = 0
$ws.Range("P23").Copy()
$ws.Range("Q23").PasteSpecial(-4122)
$ws.Range("Q23").Value = 0
$ws.Range("P24").Copy()
$ws.Range("Q24").PasteSpecial(-4122)
$ws.Range("Q24").Value = 0
$ws.Range("P25").Copy()
$ws.Range("Q25").PasteSpecial(-4122)
$ws.Range("Q25").Value = 0.1
$ws.Range("P26").Copy()
$ws.Range("Q26").PasteSpecial(-4122)
$ws.Range("Q26").Value = 0.2
$ws.Range("P27").Copy()
$ws.Range("Q27").PasteSpecial(-4122)
$ws.Range("Q27").Value = 0.1
$ws.Range("P28").Copy()
$ws.Range("Q28").PasteSpecial(-4122)
$ws.Range("Q28").Value = 0.3
$ws.Range("P29").Copy()
$ws.Range("Q29").PasteSpecial(-4122)
$ws.Range("Q29").Value = 0.4
$ws.Range("P30").Copy()
$ws.Range("Q30").PasteSpecial(-4122)
$ws.Range("Q30").Value = 0.2
$ws.Range("P31").Copy()
$ws.Range("Q31").PasteSpecial(-4122)
$ws.Range("Q31").Value = 0.2
$ws.Range("P32").Copy()
$ws.Range("Q32").PasteSpecial(-4122)
$ws.Range("Q32").Value = 0.2
$ws.Range("P33").Copy()
$ws.Range("Q33").PasteSpecial(-4122)
$ws.Range("Q33").Value = 0.1
$ws.Range("P34").Copy()
$ws.Range("Q34").PasteSpecial(-4122)
$ws.Range("P35").Copy()
$ws.Range("Q35").PasteSpecial(-4122)
$ws.Range("Q35").Value = 0
$ws.Range("P36").Copy()
$ws.Range("Q36").PasteSpecial(-4122)
$ws.Range("Q36").Value = 0.1
$ws.Range("P37").Copy()
$ws.Range("Q37").PasteSpecial(-4122)
$ws.Range("Q37").Value = 0.2

# Leave the active selection on P30, matching the saved workbook state.
$ws.Range("P30").Select()
